$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '22.392.88'
$ws.Range("E2").Value = '  +0.15%  '

# Row 3
$ws.Range("D3").Value = '1.571.01'
$ws.Range("E3").Value = '  +0.32%  '

# Row 4
$ws.Range("E4").Value = '  +0.04%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9993'
$ws.Range("E5").Value = '  -0.18%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '291.45'
$ws.Range("E6").Value = '  +0.68%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3760'
$ws.Range("E7").Value = '  +2.13%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '49.76'
$ws.Range("E8").Value = '  +0.97%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3425'
$ws.Range("E9").Value = '  +1.39%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.07631'
$ws.Range("E10").Value = '  +0.58%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.148'
$ws.Range("E11").Value = '  -1.57%  '

# Row 12
$ws.Range("E12").Value = '  +0.04%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '21.15'
$ws.Range("E13").Value = '  -0.09%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.006'
$ws.Range("E14").Value = '  -0.60%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.958'
$ws.Range("E15").Value = '  +1.14%  '

# Row 16
$ws.Range("B16").Value = 'ShibaInu'
$ws.Range("C16").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001130'
$ws.Range("E16").Value = '  -0.04%  '

# Row 17
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '1.560.03'
$ws.Range("E17").Value = '  -0.93%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '89.97'
$ws.Range("E18").Value = '  +0.88%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06718'
$ws.Range("E19").Value = '  -0.50%  '

# Row 20
$ws.Range("E20").Value = '  +0.04%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '16.75'
$ws.Range("E21").Value = '  +1.52%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.205'
$ws.Range("E22").Value = '  -0.24%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.00'
$ws.Range("E23").Value = '  +0.37%  '

# Row 24
$ws.Range("D24").Value = '22.387.94'
$ws.Range("E24").Value = '  +0.05%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.394'
$ws.Range("E25").Value = '  +0.70%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.651'
$ws.Range("E26").Value = '  -11.27%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '20.15'
$ws.Range("E27").Value = '  +1.43%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '146.72'
$ws.Range("E28").Value = '  +0.75%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.028'
$ws.Range("E29").Value = '  +1.23%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '126.50'
$ws.Range("E30").Value = '  +1.05%  '

# Row 31
$ws.Range("D31").Value = '1.741.80'
$ws.Range("E31").Value = '  -0.33%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.134'
$ws.Range("E32").Value = '  -1.85%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.994'
$ws.Range("E33").Value = '  +0.31%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.9812'
$ws.Range("E34").Value = '  -5.67%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '10.09'
$ws.Range("E35").Value = '  -1.72%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08507'
$ws.Range("E36").Value = '  +0.63%  '

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.413'
$ws.Range("E37").Value = '  +13.53%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.02541'
$ws.Range("E38").Value = '  +0.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.2313'
$ws.Range("E39").Value = '  -0.64%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.06555'
$ws.Range("E40").Value = '  +0.67%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.412'
$ws.Range("E41").Value = '  -2.34%  '

# Row 42
$ws.Range("B42").Value = 'Aptos'
$ws.Range("C42").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '11.45'
$ws.Range("E42").Value = '  -2.74%  '

# Row 43
$ws.Range("B43").Value = 'TheSandbox'
$ws.Range("C43").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6389'
$ws.Range("E43").Value = '  +0.57%  '

# Row 44
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.10'
$ws.Range("E44").Value = '  -1.65%  '

# Row 45
$ws.Range("B45").Value = 'Frax'
$ws.Range("C45").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.001'
$ws.Range("E45").Value = '  +0.06%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.795'

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5970'
$ws.Range("E47").Value = '  -0.04%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.292'
$ws.Range("E48").Value = '  +2.24%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.085'
$ws.Range("E49").Value = '  -2.02%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '124.83'
$ws.Range("E50").Value = '  +0.91%  '

# Row 51
$ws.Range("E51").Value = '  +0.76%  '
